# everyday_ver2.xlsx update: add new "2021/11/23" row to the top of each
# data table (row 2), shifting the existing history rows down by one.
#
# Sheet 1: 台指期換倉成本計算
# Sheet 2: 散戶多空力道
# Sheet 3: 三大法人買賣金額
# Sheet 4: 大盤多空點位
# Sheet 5: 期貨大額交易人未沖銷部位

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算  (A1:F5 -> A1:F6)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()

# keep the trailing blank row (previously row 5, now row 6) present so the
# sheet's used range / dimension still extends one row past the data, just
# like the original file.
$ws1.Range("F6").Borders.LineStyle = 0

$ws1.Range("A2").Value = "日期：2021/11/23"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "202201"
$ws1.Range("C2").Value = 17610
$ws1.Range("D2").Value = 4365
$ws1.Range("E2").Value = 14774790
$ws1.Range("F2").Value = 17770

# ---------------------------------------------------------------------
# Sheet 2: 散戶多空力道  (A1:B20 -> A1:B21)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()

$ws2.Range("A2").Value = "日期：2021/11/23"
$ws2.Range("B2").Value = 0.17

# ---------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額  (A1:C20 -> A1:C21)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()

$ws3.Range("A2").Value = "110年11月23日"
$ws3.Range("B2").Value = -128.1
$ws3.Range("C2").Value = -9.32

# ---------------------------------------------------------------------
# Sheet 4: 大盤多空點位  (A1:B19 -> A1:B20)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()

$ws4.Range("A2").Value = "110年11月23日"
$ws4.Range("B2").Value = 17711.53

# ---------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位  (A1:N19 -> A1:N20)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()

$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/11/23"
$ws5.Range("B2").Value = 46450
$ws5.Range("C2").Value = 51600
$ws5.Range("D2").Value = 31
$ws5.Range("E2").Value = 284
$ws5.Range("F2").Value = 25391
$ws5.Range("G2").Value = 44098
$ws5.Range("H2").Value = -728
$ws5.Range("I2").Value = -677
$ws5.Range("J2").Value = -18707
$ws5.Range("K2").Value = -51
$ws5.Range("L2").Value = 759
$ws5.Range("M2").Value = 961
$ws5.Range("N2").Value = -202
